$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert "Actual Result" column before "Result"/"Revise" ---
$ws.Range("F1").Value = "Actual Result"
$ws.Range("G1").Value = "Result"
$ws.Range("H1").Value = "Revise"

# --- Data rows: fill in "Actual Result" (column D) and refresh Expected Result (col E) ---

# Row 2
$ws.Range("C2").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D2").Value = "#jame1234"
$ws.Range("E2").Value = "เข้าสู่ระบบของผู้ดูแลระบบสำเร็จ"

# Row 3
$ws.Range("C3").Value = "Jalem4544471@mju.ac.th"
$ws.Range("D3").Value = "#jame1235"
$ws.Range("E3").Value = "กรุณากรอกข้อมูลอีเมล์ให้อยู่ในรูปแบบของ MJU เช่น (MJU6504106001@mju.ac.th)"

# Row 4
$ws.Range("C4").Value = "MJU6504106333@mju.ac.th"
$ws.Range("D4").Value = "#jame1236"
$ws.Range("E4").Value = "เข้าสู่ระบบของผู้ดูแลระบบสำเร็จ"

# Row 5
$ws.Range("C5").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D5").Value = "#jame1237"
$ws.Range("E5").Value = "กรุณากรอกข้อมูลโดยห้ามมีช่องว่าง"

# Row 6
$ws.Range("C6").Value = "mju6504106336@mju.ac.th"
$ws.Range("D6").Value = "#jame1238"
$ws.Range("E6").Value = "เข้าสู่ระบบของผู้ดูแลระบบสำเร็จ"

# Row 7
$ws.Range("C7").Value = "mju650410633@mju.ac.th"
$ws.Range("D7").Value = "#jame1239"
$ws.Range("E7").Value = "กรุณากรอกความยาวให้มีขนาด 13 ตัวพอดี"

# Row 8
$ws.Range("C8").Value = "mju65041063361@mju.ac.th"
$ws.Range("D8").Value = "#jame1240"
$ws.Range("E8").Value = "กรุณากรอกความยาวให้มีขนาด 13 ตัวพอดี"

# Row 9 (C9 stays blank)
$ws.Range("D9").Value = "#jame1241"
$ws.Range("E9").Value = "กรุณากรอกอีเมล์รูปแบบของอีเมล MJU"

# Row 10
$ws.Range("C10").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D10").Value = "#jame1234"
$ws.Range("E10").Value = "เข้าสู่ระบบของผู้ดูแลระบบสำเร็จ"

# Row 11
$ws.Range("C11").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D11").Value = "#เจม1234"
$ws.Range("E11").Value = "กรุณากรอกข้อมูลเป็นตัวอักษรภาษาอังกฤษหรือตัวเลขรวมอักขระพิเศษ"

# Row 12
$ws.Range("C12").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D12").Value = "#jame_12"
$ws.Range("E12").Value = "เข้าสู่ระบบของผู้ดูแลระบบสำเร็จ"

# Row 13
$ws.Range("C13").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D13").Value = "#jame_123"
$ws.Range("E13").Value = "เข้าสู่ระบบของผู้ดูแลระบบสำเร็จ"

# Row 14
$ws.Range("C14").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D14").Value = "#jame_123456789"
$ws.Range("E14").Value = "เข้าสู่ระบบของผู้ดูแลระบบสำเร็จ"

# Row 15
$ws.Range("C15").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D15").Value = "#jame12345678910"
$ws.Range("E15").Value = "เข้าสู่ระบบของผู้ดูแลระบบสำเร็จ"

# Row 16
$ws.Range("C16").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D16").Value = "#jame12"
$ws.Range("E16").Value = "กรุณากรอกข้อมูลได้ตั้งแต่ 8 ถึง 16 ตัวอักษร"

# Row 17
$ws.Range("C17").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D17").Value = "##jame123456789__"
$ws.Range("E17").Value = "กรุณากรอกข้อมูลได้ตั้งแต่ 8 ถึง 16 ตัวอักษร"

# Row 18
$ws.Range("C18").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D18").Value = "#jame1234"
$ws.Range("E18").Value = "เข้าสู่ระบบของผู้ดูแลระบบสำเร็จ"

# Row 19
$ws.Range("C19").Value = "MJU6504106336@mju.ac.th"
$ws.Range("D19").Value = "#jame 1234"
$ws.Range("E19").Value = "กรุณากรอกข้อมูลห้ามมีช่องว่างระหว่าตัวอักษร"

# Row 20 (D20 stays blank)
$ws.Range("C20").Value = "MJU6504106336@mju.ac.th"
$ws.Range("E20").Value = "กรุณากรอกรหัสผ่าน"

# --- Update selection/view state ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("E25").Select()
